$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set Runmode column (C) to "Y" for rows 3-5 so that every test case runs
# (previously only Suite A's first case ran, others were "N"/skipped).
$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"
